$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain text (matches original formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = "67.490.07"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "2.628.84"
$ws.Range("E3").Value = "  -1.80%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "594.67"
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("D6").Value = "168.43"
$ws.Range("E6").Value = "  +0.99%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("D9").Value = "2.628.56"
$ws.Range("E9").Value = "  -1.76%  "

$ws.Range("E10").Value = "  -2.11%  "

$ws.Range("E11").Value = "  +1.22%  "

$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("D13").Value = "5.23"
$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("D14").Value = "27.63"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").Value = "3.108.61"
$ws.Range("E15").Value = "  -1.87%  "

$ws.Range("E16").Value = "  -1.49%  "

$ws.Range("D17").Value = "67.349.42"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").Value = "2.628.23"
$ws.Range("E18").Value = "  -1.88%  "

$ws.Range("D19").Value = "11.99"
$ws.Range("E19").Value = "  +2.01%  "

$ws.Range("D20").Value = "8.05"
$ws.Range("E20").Value = "  +4.27%  "

$ws.Range("D21").Value = "357.26"
$ws.Range("E21").Value = "  -1.91%  "

$ws.Range("D22").Value = "4.31"
$ws.Range("E22").Value = "  -1.76%  "

$ws.Range("E23").Value = "  -3.25%  "

$ws.Range("E24").Value = "  -4.25%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").Value = "10.33"
$ws.Range("E26").Value = "  +2.48%  "

$ws.Range("D27").Value = "69.77"
$ws.Range("E27").Value = "  -1.57%  "

$ws.Range("D28").Value = "2.760.06"
$ws.Range("E28").Value = "  -2.35%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -1.72%  "

$ws.Range("D31").Value = "547.32"
$ws.Range("E31").Value = "  -1.96%  "

$ws.Range("E32").Value = "  -1.11%  "

$ws.Range("D33").Value = "1.36"
$ws.Range("E33").Value = "  -2.98%  "

$ws.Range("E34").Value = "  -2.04%  "

$ws.Range("E35").Value = "  +4.04%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Value = "1.50"
$ws.Range("E37").Value = "  -4.08%  "

$ws.Range("D38").Value = "156.91"
$ws.Range("E38").Value = "  +0.48%  "

$ws.Range("D39").Value = "19.02"
$ws.Range("E39").Value = "  -2.63%  "

$ws.Range("E40").Value = "  -2.16%  "

$ws.Range("E41").Value = "  -0.78%  "

$ws.Range("E42").Value = "  +1.90%  "

$ws.Range("D43").Value = "5.22"
$ws.Range("E43").Value = "  -2.05%  "

$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("E45").Value = "  -3.79%  "

$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("D47").Value = "152.94"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("E48").Value = "  -1.95%  "

$ws.Range("E49").Value = "  -1.53%  "

$ws.Range("E50").Value = "  -1.35%  "

$ws.Range("E51").Value = "  -1.07%  "
